$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update 想去人数 (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5390
$wsExpo.Range("F3").Value = 586
$wsExpo.Range("F4").Value = 11447
$wsExpo.Range("F5").Value = 281
$wsExpo.Range("F8").Value = 255
$wsExpo.Range("F9").Value = 986

# Sheet "全部类型" (all types) - same underlying rows, duplicated data
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5390
$wsAll.Range("F5").Value = 586
$wsAll.Range("F7").Value = 11447
$wsAll.Range("F8").Value = 281
$wsAll.Range("F13").Value = 255
$wsAll.Range("F14").Value = 986
